# The deck's single theme ("Integral") is swapped back to the stock
# PowerPoint "Office Theme" colour palette. dk1/lt1 (black/white) are
# unchanged; the other ten theme colours (dk2, lt2, accent1-6, hlink,
# folHlink) are updated to the default Office theme RGB values via the
# presentation's ThemeColorScheme (index order matches the OOXML
# <a:clrScheme> child order: dk1, lt1, dk2, lt2, accent1..6, hlink,
# folHlink).

$p = $ppt.ActivePresentation

# msoThemeColor-style index -> target Office-theme RGB (as r + g*256 + b*65536,
# matching the VBA RGB() long packing used by ThemeColor.RGB).
$officeColors = @{
    3  = (0x44 + (0x54 * 256) + (0x6A * 65536))   # dk2      -> 44546A
    4  = (0xE7 + (0xE6 * 256) + (0xE6 * 65536))   # lt2      -> E7E6E6
    5  = (0x5B + (0x9B * 256) + (0xD5 * 65536))   # accent1  -> 5B9BD5
    6  = (0xED + (0x7D * 256) + (0x31 * 65536))   # accent2  -> ED7D31
    7  = (0xA5 + (0xA5 * 256) + (0xA5 * 65536))   # accent3  -> A5A5A5
    8  = (0xFF + (0xC0 * 256) + (0x00 * 65536))   # accent4  -> FFC000
    9  = (0x44 + (0x72 * 256) + (0xC4 * 65536))   # accent5  -> 4472C4
    10 = (0x70 + (0xAD * 256) + (0x47 * 65536))   # accent6  -> 70AD47
    11 = (0x05 + (0x63 * 256) + (0xC1 * 65536))   # hlink    -> 0563C1
    12 = (0x95 + (0x4F * 256) + (0x72 * 65536))   # folHlink -> 954F72
}

$tcs = $p.Slides.Item(1).ThemeColorScheme
foreach ($idx in $officeColors.Keys) {
    $tcs.Colors($idx).RGB = $officeColors[$idx]
}
